$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value reads as a plain number need to be
# pre-formatted as Text so Excel keeps them as literal strings, just
# like the rest of the (text-typed) Price column.
$textCells = @("D5","D10","D14","D15","D16","D19","D21","D22","D24","D25","D32","D34","D38","D45","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.188.10"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.643.28"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "217.19"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "20.04"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "1.873.09"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.631.43"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "4.17"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("D16").Value = "67.37"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "27.144.49"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "219.09"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +2.77%  "
$ws.Range("D22").Value = "2.58"
$ws.Range("E22").Value = "  +6.62%  "
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "9.22"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "147.63"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").Value = "1.272.74"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "0.859"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +7.23%  "
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").Value = "1.783.32"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "61.90"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "0.0977"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "7.64"
$ws.Range("E51").Value = "  +1.32%  "
